$wb = $excel.ActiveWorkbook

# The workbook has duplicated data across the "展览" (Exhibitions) and
# "全部类型" (All Types) sheets. Row 4 is "南宁·2024良牙动漫秋季盛典（秋典）"
# and row 5 is "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini". Column F holds the
# "想去人数" (number of people interested) counter, which was updated on a
# refresh of the source data.

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 3729
    $ws.Range("F5").Value = 385
}
